# Add a new "region" (short / English-friendly) column to the Chile regions
# sheet. The historical column B ("region") held the long official Spanish
# region name (e.g. "I Región de Tarapacá") - that column is renamed to
# "region_long_name" and a brand-new column C ("region") is inserted with the
# short region names. Everything to the right (region_code, region_number,
# population) shifts one column over, from C:E to D:F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at C - shifts former C/D/E (region_code,
# region_number, population) one slot right to D/E/F, leaves A/B untouched.
$ws.Columns("C:C").Insert()

# Header row.
$ws.Range("B1").Value = "region_long_name"
$ws.Range("C1").Value = "region"

# New short region names for column C, row by row.
$ws.Range("C2").Value = "Chile"
$ws.Range("C3").Value = "Tarapacá"
$ws.Range("C4").Value = "Antofagasta"
$ws.Range("C5").Value = "Atacama"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("C7").Value = "Valparaíso"
$ws.Range("C8").Value = "Libertador General Bernardo OHiggins"
$ws.Range("C9").Value = "Maule"
$ws.Range("C10").Value = "Biobio"
$ws.Range("C11").Value = "Ñuble"
$ws.Range("C12").Value = "Araucanía"
$ws.Range("C13").Value = "Los Lagos"
$ws.Range("C14").Value = "General Carlos Ibáñez del Campo's Aysén"
$ws.Range("C15").Value = "Magellan and Chilean Antarctica"
$ws.Range("C16").Value = "Santiago Metropolitan"
$ws.Range("C17").Value = "Los Ríos"
$ws.Range("C18").Value = "Arica y Parinacota"

# Match the authored column width/formatting for the new column. Excel
# stores column width in "character" units derived from pixels (padding +
# max digit width), so asking for an even 48 characters round-trips to
# 48.83.. on disk; back the desired value off by the padding (5/6 of a
# character here) so the persisted <col width="..."/> comes out at exactly
# 48, matching the authored file.
$ws.Columns("C:C").ColumnWidth = 47.1666667

# The sheet-scoped "_FilterDatabase" name (behind the hidden filter
# dropdowns) still points at the old A1:E18 block; widen it to cover the
# newly-inserted column.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$F`$18"
    }
}

# Selection moved to B10 in the saved workbook.
$ws.Range("B10").Select()
